$d = $word.ActiveDocument
$q = [char]34

# Make sure literal straight quotes are inserted (avoid "smart quote" autocorrect)
try { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}
try { $word.Options.AutoFormatReplaceQuotes = $false } catch {}
try { $word.AutoCorrect.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}

function Replace-ExactOnce([string]$findText, [string]$replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $replaceText
    }
    return $found
}

# --- "Username: ?" / "Password: ?" / "Email: ?" -> wrap the "?" in straight quotes ---
# login block
Replace-ExactOnce "Username: ?" ("Username: " + $q + "?" + $q) | Out-Null
Replace-ExactOnce "Password: ?" ("Password: " + $q + "?" + $q) | Out-Null
# signup block
Replace-ExactOnce "Username: ?" ("Username: " + $q + "?" + $q) | Out-Null
Replace-ExactOnce "Password: ?" ("Password: " + $q + "?" + $q) | Out-Null
Replace-ExactOnce "Email: ?" ("Email: " + $q + "?" + $q) | Out-Null

# --- JoinRoom request fields: only RoomName and AnswerTime get wrapped in quotes ---
Replace-ExactOnce "RoomName:?" ("RoomName:" + $q + "?" + $q) | Out-Null
Replace-ExactOnce "AnswerTime:?" ("AnswerTime:" + $q + "?" + $q) | Out-Null

# --- remove the space before the opening bracket ---
Replace-ExactOnce "Rooms: [" "Rooms[" | Out-Null
Replace-ExactOnce "Names: [" "Names[" | Out-Null

# --- "Highscores :[" -> "Highscores[" (drop the trailing space + colon) ---
Replace-ExactOnce "Highscores :[" "Highscores[" | Out-Null
